$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.021112519215942
$ws.Range("D2").Value = 1.040651558982093
$ws.Range("E2").Value = 1.035020050190573
$ws.Range("F2").Value = 1.046557582252405
$ws.Range("I2").Value = 1.058442608661789
$ws.Range("J2").Value = 1.042679587686169
$ws.Range("K2").Value = 1.051572178504137
$ws.Range("L2").Value = 1.04601214233958
$ws.Range("M2").Value = 1.057404274665112
$ws.Range("N2").Value = 1.018039286016822
$ws.Range("R2").Value = 1.047536224427962
$ws.Range("C3").Value = 1.024486619651238
$ws.Range("D3").Value = 1.043003737191843
$ws.Range("E3").Value = 1.037707545082601
$ws.Range("F3").Value = 1.049378277891512
$ws.Range("I3").Value = 1.059397942490786
$ws.Range("J3").Value = 1.044331305128987
$ws.Range("K3").Value = 1.053114427441856
$ws.Range("L3").Value = 1.047879655006826
$ws.Range("M3").Value = 1.059416056866395
$ws.Range("N3").Value = 1.018592509495565
$ws.Range("R3").Value = 1.048624069083575
$ws.Range("C4").Value = 1.026635590434357
$ws.Range("D4").Value = 1.044505536973646
$ws.Range("E4").Value = 1.039424675790824
$ws.Range("F4").Value = 1.051180032438341
$ws.Range("I4").Value = 1.059998164464554
$ws.Range("J4").Value = 1.045381652710731
$ws.Range("K4").Value = 1.054094928591095
$ws.Range("L4").Value = 1.049069471280263
$ws.Range("M4").Value = 1.060697639418078
$ws.Range("N4").Value = 1.018944383108633
$ws.Range("Q4").Value = 1.019999999999999
$ws.Range("R4").Value = 1.049318235139487
$ws.Range("C5").Value = 1.027533817103148
$ws.Range("D5").Value = 1.045135967811302
$ws.Range("E5").Value = 1.040144201080794
$ws.Range("F5").Value = 1.051935185618755
$ws.Range("I5").Value = 1.060248442793159
$ws.Range("J5").Value = 1.045821628257418
$ws.Range("K5").Value = 1.054506804326537
$ws.Range("L5").Value = 1.049568034580036
$ws.Range("M5").Value = 1.061234823761097
$ws.Range("N5").Value = 1.019092069002579
$ws.Range("R5").Value = 1.049616537535769
$ws.Range("C6").Value = 1.027687568469575
$ws.Range("D6").Value = 1.045246174060326
$ws.Range("E6").Value = 1.040268060206498
$ws.Range("F6").Value = 1.052065514204104
$ws.Range("I6").Value = 1.060292871420458
$ws.Range("J6").Value = 1.045898571345183
$ws.Range("K6").Value = 1.054580302505954
$ws.Range("L6").Value = 1.049654800908874
$ws.Range("M6").Value = 1.061328557452313
$ws.Range("N6").Value = 1.019118231112644
$ws.Range("R6").Value = 1.049677138202888
$ws.Range("C7").Value = 1.02665674443439
$ws.Range("D7").Value = 1.044526376073138
$ws.Range("E7").Value = 1.039443284800421
$ws.Range("F7").Value = 1.051200470055424
$ws.Range("I7").Value = 1.060008630033672
$ws.Range("J7").Value = 1.04539645196512
$ws.Range("K7").Value = 1.05411271053212
$ws.Range("L7").Value = 1.049085026451909
$ws.Range("M7").Value = 1.060715051480283
$ws.Range("N7").Value = 1.018950256380443
$ws.Range("R7").Value = 1.049350664181052
$ws.Range("C8").Value = 1.022271371244295
$ws.Range("D8").Value = 1.041466013788186
$ws.Range("E8").Value = 1.035944047910507
$ws.Range("F8").Value = 1.04752854115842
$ws.Range("I8").Value = 1.058777973944969
$ws.Range("J8").Value = 1.043252679224234
$ws.Range("K8").Value = 1.052112135295665
$ws.Range("L8").Value = 1.046658188530972
$ws.Range("M8").Value = 1.05810101040791
$ws.Range("N8").Value = 1.018232395689428
$ws.Range("R8").Value = 1.04794063689029
$ws.Range("C9").Value = 1.014256361601758
$ws.Range("D9").Value = 1.035893714398888
$ws.Range("E9").Value = 1.029590213260965
$ws.Range("F9").Value = 1.040855961493183
$ws.Range("I9").Value = 1.056457919193249
$ws.Range("J9").Value = 1.039315543836891
$ws.Range("K9").Value = 1.048430097842576
$ws.Range("L9").Value = 1.042220782937342
$ws.Range("M9").Value = 1.053319218252912
$ws.Range("N9").Value = 1.016913108540349
$ws.Range("R9").Value = 1.04533411531074
$ws.Range("C10").Value = 1.008717564172759
$ws.Range("D10").Value = 1.03206758459872
$ws.Range("E10").Value = 1.025231815894301
$ws.Range("F10").Value = 1.036276524785798
$ws.Range("I10").Value = 1.054814914719075
$ws.Range("J10").Value = 1.036589526853157
$ws.Range("K10").Value = 1.045881020951913
$ws.Range("L10").Value = 1.039159971576085
$ws.Range("M10").Value = 1.050020294286998
$ws.Range("N10").Value = 1.016000549195791
$ws.Range("R10").Value = 1.043548421987619
$ws.Range("C11").Value = 1.006279997192091
$ws.Range("D11").Value = 1.030396600902873
$ws.Range("E11").Value = 1.023323936539085
$ws.Range("F11").Value = 1.034272074548058
$ws.Range("I11").Value = 1.054087627993862
$ws.Range("J11").Value = 1.035393834680665
$ws.Range("K11").Value = 1.044767179553239
$ws.Range("L11").Value = 1.037819006063692
$ws.Range("M11").Value = 1.048575372327864
$ws.Range("N11").Value = 1.015601761618507
$ws.Range("R11").Value = 1.042793397446252
$ws.Range("C12").Value = 1.005360413633655
$ws.Range("D12").Value = 1.02976317091629
$ws.Range("E12").Value = 1.022604163795796
$ws.Range("F12").Value = 1.033515288621654
$ws.Range("I12").Value = 1.0538090175033
$ws.Range("J12").Value = 1.034939631045417
$ws.Range("K12").Value = 1.044341538465608
$ws.Range("L12").Value = 1.037310749293661
$ws.Range("M12").Value = 1.048027375213973
$ws.Range("N12").Value = 1.015449572619041
$ws.Range("R12").Value = 1.042492461185037
$ws.Range("C13").Value = 1.005556795034378
$ws.Range("D13").Value = 1.029897632094313
$ws.Range("E13").Value = 1.022757571580052
$ws.Range("F13").Value = 1.033676498763897
$ws.Range("I13").Value = 1.053868023417355
$ws.Range("J13").Value = 1.035036060915086
$ws.Range("K13").Value = 1.044431424367409
$ws.Range("L13").Value = 1.037418769514938
$ws.Range("M13").Value = 1.04814378456716
$ws.Range("N13").Value = 1.015481738530812
$ws.Range("R13").Value = 1.042553544521486
$ws.Range("C14").Value = 1.006203619862927
$ws.Range("D14").Value = 1.030343642621753
$ws.Range("E14").Value = 1.023264024058701
$ws.Range("F14").Value = 1.034209044047646
$ws.Range("I14").Value = 1.054064278294114
$ws.Range("J14").Value = 1.035355867768115
$ws.Range("K14").Value = 1.044731394668426
$ws.Range("L14").Value = 1.037776569477964
$ws.Range("M14").Value = 1.048529593048108
$ws.Range("N14").Value = 1.01558897853598
$ws.Range("R14").Value = 1.042767039731368
$ws.Range("C15").Value = 1.006603543005218
$ws.Range("D15").Value = 1.030621058611624
$ws.Range("E15").Value = 1.023577810882443
$ws.Range("F15").Value = 1.03453916746939
$ws.Range("I15").Value = 1.054186542908613
$ws.Range("J15").Value = 1.035554721055513
$ws.Range("K15").Value = 1.044918869315608
$ws.Range("L15").Value = 1.037998836238802
$ws.Range("M15").Value = 1.048769373490227
$ws.Range("N15").Value = 1.015655945537332
$ws.Range("R15").Value = 1.042905400299189
$ws.Range("B16").Value = 1.05
$ws.Range("C16").Value = 1.008899413525869
$ws.Range("D16").Value = 1.032205740832796
$ws.Range("E16").Value = 1.025378548768488
$ws.Range("F16").Value = 1.036432310636098
$ws.Range("I16").Value = 1.054878645227476
$ws.Range("J16").Value = 1.036688667867439
$ws.Range("K16").Value = 1.045981834534268
$ws.Range("L16").Value = 1.039268842998598
$ws.Range("M16").Value = 1.05013867484995
$ws.Range("N16").Value = 1.016036062860301
$ws.Range("R16").Value = 1.043660074769522
$ws.Range("C17").Value = 1.010325660364207
$ws.Range("D17").Value = 1.033191715811109
$ws.Range("E17").Value = 1.026499331740858
$ws.Range("F17").Value = 1.037610436435298
$ws.Range("I17").Value = 1.0553061079178
$ws.Range("J17").Value = 1.037392741047175
$ws.Range("K17").Value = 1.046641769635202
$ws.Range("L17").Value = 1.040058154855817
$ws.Range("M17").Value = 1.050989659799561
$ws.Range("N17").Value = 1.016272105825286
$ws.Range("R17").Value = 1.044129224344634
$ws.Range("C18").Value = 1.011146148609521
$ws.Range("D18").Value = 1.033755354106126
$ws.Range("E18").Value = 1.027143561872484
$ws.Range("F18").Value = 1.03828705129673
$ws.Range("I18").Value = 1.055548126683016
$ws.Range("J18").Value = 1.037794581960602
$ws.Range("K18").Value = 1.047015798595175
$ws.Range("L18").Value = 1.040509652213452
$ws.Range("M18").Value = 1.051476080626053
$ws.Range("N18").Value = 1.01640609872046
$ws.Range("R18").Value = 1.044382154910104
$ws.Range("C19").Value = 1.011430125210455
$ws.Range("D19").Value = 1.033953717068165
$ws.Range("E19").Value = 1.027367655274973
$ws.Range("F19").Value = 1.038522800721964
$ws.Range("I19").Value = 1.055634063175525
$ws.Range("J19").Value = 1.037936017934977
$ws.Range("K19").Value = 1.047149475204452
$ws.Range("L19").Value = 1.040668033089114
$ws.Range("M19").Value = 1.051646975495264
$ws.Range("N19").Value = 1.016453837746857
$ws.Range("R19").Value = 1.044482991832165
$ws.Range("C20").Value = 1.010172641055284
$ws.Range("D20").Value = 1.033085557761192
$ws.Range("E20").Value = 1.026378907039787
$ws.Range("F20").Value = 1.037483817233451
$ws.Range("I20").Value = 1.055260108624824
$ws.Range("J20").Value = 1.037316979167367
$ws.Range("K20").Value = 1.046570558881803
$ws.Range("L20").Value = 1.03997324943448
$ws.Range("M20").Value = 1.05089809581762
$ws.Range("N20").Value = 1.016246648032507
$ws.Range("R20").Value = 1.04407755536819
$ws.Range("C21").Value = 1.006020320141333
$ws.Range("D21").Value = 1.030221545685638
$ws.Range("E21").Value = 1.023121853524208
$ws.Range("F21").Value = 1.034060063070186
$ws.Range("I21").Value = 1.054011815361583
$ws.Range("J21").Value = 1.03526847303042
$ws.Range("K21").Value = 1.044652123324355
$ws.Range("L21").Value = 1.037678028066215
$ws.Range("M21").Value = 1.048423665201915
$ws.Range("N21").Value = 1.015560476636563
$ws.Range("R21").Value = 1.042724369072539
$ws.Range("C22").Value = 1.003354491623341
$ws.Range("D22").Value = 1.028382617431183
$ws.Range("E22").Value = 1.021036230728778
$ws.Range("F22").Value = 1.031866485083089
$ws.Range("I22").Value = 1.053198385007632
$ws.Range("J22").Value = 1.033948155362061
$ws.Range("K22").Value = 1.043412026028844
$ws.Range("L22").Value = 1.036202165542539
$ws.Range("M22").Value = 1.0468320335557
$ws.Range("N22").Value = 1.015117296855181
$ws.Range("R22").Value = 1.041834219179505
$ws.Range("C23").Value = 1.00476388757782
$ws.Range("D23").Value = 1.029349173234784
$ws.Range("E23").Value = 1.02213652671592
$ws.Range("F23").Value = 1.033023182812535
$ws.Range("I23").Value = 1.053625477718907
$ws.Range("J23").Value = 1.034642431958401
$ws.Range("K23").Value = 1.04406092193057
$ws.Range("L23").Value = 1.036978887560447
$ws.Range("M23").Value = 1.047669309654603
$ws.Range("N23").Value = 1.01534936564245
$ws.Range("R23").Value = 1.042283472126986
$ws.Range("C24").Value = 1.010228092246846
$ws.Range("D24").Value = 1.033115434058099
$ws.Range("E24").Value = 1.026419859893316
$ws.Range("F24").Value = 1.037525807325685
$ws.Range("I24").Value = 1.05527052612605
$ws.Range("J24").Value = 1.03733802649823
$ws.Range("K24").Value = 1.046584934914745
$ws.Range("L24").Value = 1.039998372255256
$ws.Range("M24").Value = 1.050924487675102
$ws.Range("N24").Value = 1.016252179638018
$ws.Range("R24").Value = 1.044060612967343
$ws.Range("C25").Value = 1.016376639720076
$ws.Range("D25").Value = 1.037373180200198
$ws.Range("E25").Value = 1.03126840261744
$ws.Range("F25").Value = 1.042620025073181
$ws.Range("I25").Value = 1.057085743463552
$ws.Range("J25").Value = 1.04036516027235
$ws.Range("K25").Value = 1.049417761713629
$ws.Range("L25").Value = 1.04339996132506
$ws.Range("M25").Value = 1.05459088018758
$ws.Range("N25").Value = 1.017266296642117
$ws.Range("R25").Value = 1.046060693874175
